$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CS2025_00112"
$ws.Range("B2").Value = "Alice"
$ws.Range("C2").Value = " alice.johnson@example.com"
$ws.Range("D2").Value = " Technical"
$ws.Range("E2").Value = " Unable to log in to my account after recent update."

$ws.Range("A3").Value = "CS2025_00113"
$ws.Range("B3").Value = "Amanda"
$ws.Range("C3").Value = " amanda.king@example.com"
$ws.Range("D3").Value = " Technical"
$ws.Range("E3").Value = " I am getting frequent error messages on the desktop application."

$ws.Range("A4").Value = "CS2025_00114"
$ws.Range("B4").Value = "Benjamin"
$ws.Range("C4").Value = " benjamin.lewis@example.com"
$ws.Range("D4").Value = " Billing"
$ws.Range("E4").Value = " I was charged for a subscription I canceled last month."

$ws.Range("A5").Value = "CS2025_00115"
$ws.Range("B5").Value = "Bob"
$ws.Range("C5").Value = " bob.smith@example.com"
$ws.Range("D5").Value = " Billing"
$ws.Range("E5").Value = " Incorrect charge appeared on my credit card statement."

$ws.Range("A6").Value = "CS2025_00116"
$ws.Range("B6").Value = "Catherine"
$ws.Range("C6").Value = " catherine.martin@example.com"
$ws.Range("D6").Value = " General"
$ws.Range("E6").Value = " What security features are included in your software?"

$ws.Range("A7").Value = "CS2025_00117"
$ws.Range("B7").Value = "Clara"
$ws.Range("C7").Value = " clara.martinez@example.com"
$ws.Range("D7").Value = " General"
$ws.Range("E7").Value = " Do you offer discounts for bulk purchases?"

$ws.Range("A8").Value = "CS2025_00118"
$ws.Range("B8").Value = "Daniel"
$ws.Range("C8").Value = " daniel.turner@example.com"
$ws.Range("D8").Value = " Technical"
$ws.Range("E8").Value = " The system doesn't recognize my external hard drive."

$ws.Range("A9").Value = "CS2025_00119"
$ws.Range("B9").Value = "David"
$ws.Range("C9").Value = " david.lee@example.com"
$ws.Range("D9").Value = " Technical"
$ws.Range("E9").Value = " Website shows a 404 error when trying to access the support page."

$ws.Range("A10").Value = "CS2025_00120"
$ws.Range("B10").Value = "Emma"
$ws.Range("C10").Value = " emma.brown@example.com"
$ws.Range("D10").Value = " Billing"
$ws.Range("E10").Value = " I need a copy of my last three invoices for reimbursement purposes."

$ws.Range("A11").Value = "CS2025_00121"
$ws.Range("B11").Value = "Frank"
$ws.Range("C11").Value = " frank.harris@example.com"
$ws.Range("D11").Value = " General"
$ws.Range("E11").Value = " What are your customer service operating hours?"

$ws.Range("A12").Value = "CS2025_00122"
$ws.Range("B12").Value = "Grace"
$ws.Range("C12").Value = " grace.wilson@example.com"
$ws.Range("D12").Value = " Technical"
$ws.Range("E12").Value = " My mobile app keeps crashing whenever I try to open the settings page."

$ws.Range("A13").Value = "CS2025_00123"
$ws.Range("B13").Value = "Harry"
$ws.Range("C13").Value = " harry.adams@example.com"
$ws.Range("D13").Value = " Billing"
$ws.Range("E13").Value = " I accidentally made a duplicate payment. Can I get a refund?"

$ws.Range("A14").Value = "CS2025_00124"
$ws.Range("B14").Value = "Isabella"
$ws.Range("C14").Value = " isabella.white@example.com"
$ws.Range("D14").Value = " General"
$ws.Range("E14").Value = " Can I change my subscription plan mid-cycle without penalties?"

$ws.Range("A15").Value = "CS2025_00125"
$ws.Range("B15").Value = "Jack"
$ws.Range("C15").Value = " jack.davis@example.com"
$ws.Range("D15").Value = " Technical"
$ws.Range("E15").Value = " I am experiencing connectivity issues with your cloud service."

$ws.Range("A16").Value = "CS2025_00126"
$ws.Range("B16").Value = "Karen"
$ws.Range("C16").Value = " karen.mitchell@example.com"
$ws.Range("D16").Value = " Billing"
$ws.Range("E16").Value = " My subscription renewal failed, and my account is deactivated."

$ws.Range("A17").Value = "CS2025_00127"
$ws.Range("B17").Value = "Liam"
$ws.Range("C17").Value = " liam.robinson@example.com"
$ws.Range("D17").Value = " Technical"
$ws.Range("E17").Value = " The software update failed to install on my device."

$ws.Range("A18").Value = "CS2025_00128"
$ws.Range("B18").Value = "Mia"
$ws.Range("C18").Value = " mia.scott@example.com"
$ws.Range("D18").Value = " General"
$ws.Range("E18").Value = " Do you have a physical store near my location?"

$ws.Range("A19").Value = "CS2025_00129"
$ws.Range("B19").Value = "Noah"
$ws.Range("C19").Value = " noah.walker@example.com"
$ws.Range("D19").Value = " Technical"
$ws.Range("E19").Value = " I need help configuring my email on a new device."

$ws.Range("A20").Value = "CS2025_00130"
$ws.Range("B20").Value = "Olivia"
$ws.Range("C20").Value = " olivia.young@example.com"
$ws.Range("D20").Value = " Billing"
$ws.Range("E20").Value = " Can I split my annual payment into monthly installments?"

$ws.Range("A21").Value = "CS2025_00131"
$ws.Range("B21").Value = "Paul"
$ws.Range("C21").Value = " paul.baker@example.com"
$ws.Range("D21").Value = " Technical"
$ws.Range("E21").Value = " I am unable to reset my password due to an error."

$ws.Range("A22").Value = "CS2025_00132"
$ws.Range("B22").Value = "Quinn"
$ws.Range("C22").Value = " quinn.moore@example.com"
$ws.Range("D22").Value = " General"
$ws.Range("E22").Value = " Can you explain the differences between your plans?"

$ws.Range("A23").Value = "CS2025_00133"
$ws.Range("B23").Value = "Rachel"
$ws.Range("C23").Value = " rachel.thomas@example.com"
$ws.Range("D23").Value = " Billing"
$ws.Range("E23").Value = " Why was I charged twice for the same transaction?"

$ws.Range("A24").Value = "CS2025_00134"
$ws.Range("B24").Value = "Sophia"
$ws.Range("C24").Value = " sophia.hernandez@example.com"
$ws.Range("D24").Value = " Technical"
$ws.Range("E24").Value = " My printer isn't connecting to the Wi-Fi network."

$ws.Range("A25").Value = "CS2025_00135"
$ws.Range("B25").Value = "Tom"
$ws.Range("C25").Value = " tom.wilson@example.com"
$ws.Range("D25").Value = " General"
$ws.Range("E25").Value = " Is there a user manual for your latest product?"

$ws.Range("A26").Value = "CS2025_00136"
$ws.Range("B26").Value = "Uma"
$ws.Range("C26").Value = " uma.carter@example.com"
$ws.Range("D26").Value = " Technical"
$ws.Range("E26").Value = " The application freezes when I try to upload large files."

$ws.Range("A27").Value = "CS2025_00137"
$ws.Range("B27").Value = "Victor"
$ws.Range("C27").Value = " victor.evans@example.com"
$ws.Range("D27").Value = " Billing"
$ws.Range("E27").Value = " I need assistance updating my payment method for my account."

$ws.Range("A28").Value = "CS2025_00138"
$ws.Range("B28").Value = "Wendy"
$ws.Range("C28").Value = " wendy.green@example.com"
$ws.Range("D28").Value = " General"
$ws.Range("E28").Value = " Do you provide training for your enterprise software?"

$ws.Range("A29").Value = "CS2025_00139"
$ws.Range("B29").Value = "Xavier"
$ws.Range("C29").Value = " xavier.hill@example.com"
$ws.Range("D29").Value = " Technical"
$ws.Range("E29").Value = " I can't access my account because of a two-factor authentication issue."

$ws.Range("A30").Value = "CS2025_00140"
$ws.Range("B30").Value = "Yolanda"
$ws.Range("C30").Value = " yolanda.hughes@example.com"
$ws.Range("D30").Value = " Billing"
$ws.Range("E30").Value = " How can I get a detailed breakdown of my charges for the last six months?"

$ws.Range("A31").Value = "CS2025_00141"
$ws.Range("B31").Value = "Zachary"
$ws.Range("C31").Value = " zachary.perez@example.com"
$ws.Range("D31").Value = " General"
$ws.Range("E31").Value = " Does your service offer multi - language support?"

